# Adding new read data from excelsheet and pass it into datadriver
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("TestData1")
$ws2 = $wb.Worksheets.Item("sample")

# --- Rename second sheet ("sample" -> "TestData2") ---
$ws2.Name = "TestData2"

# --- Populate the new TestData2 sheet with login/result data ---
# Fill columns A & B first (row by row), then column C, so the shared-string
# table is built up in the same order as the authored workbook.
$ws2.Range("A1").Value = "username"
$ws2.Range("B1").Value = "password"
$ws2.Range("A2").Value = "dexcomnew98@gmail.com"
$ws2.Range("B2").Value = "User"
$ws2.Range("A3").Value = "user@gmail.com"
$ws2.Range("B3").Value = "User123#"
$ws2.Range("A4").Value = "user@gmail.com"
$ws2.Range("B4").Value = "User123#"

$ws2.Range("C1").Value = "message"
$ws2.Range("C2").Value = "Incorrect email or password."
$ws2.Range("C3").Value = "Incorrect email or password."
$ws2.Range("C4").Value = "Incorrect"

# Highlight the header row
$ws2.Range("A1:C1").Interior.Color = 65535

# Turn the email entries in column A into mailto hyperlinks
$ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:dexcomnew98@gmail.com")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "mailto:user@gmail.com")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "mailto:user@gmail.com")

# Size the new columns to fit their content
$ws2.Columns("A:C").AutoFit()

# --- sheet1 cosmetic tweaks: widen column B, move selection ---
$ws1.Columns.Item(2).AutoFit()
[void]$ws1.Range("D16").Select()

# --- Make TestData2 the active/selected sheet & cell ---
$ws2.Activate()
[void]$ws2.Range("L11").Select()
